$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = "items"
$ws.Range("C1").Value = "path"

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "digital_elisa_protocol"
$ws.Range("C2").Value = "../files/info/digital_elisa_protocol.xlsx"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "beads"
$ws.Range("C3").Value = "../files/info/beads.xlsx"

# Styles: bold font, thin border around, centered horizontally, top vertically.
# Applying per single-cell keeps the generated style table compact (one new xf).
foreach ($addr in @("B1", "C1", "A2", "A3")) {
    $r = $ws.Range($addr)
    $f = $r.Font
    $f.Bold = $true
    $b = $r.Borders
    $b.LineStyle = 1
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4160
}

$ws.Range("H14").Select()
